$wb = $excel.ActiveWorkbook

# Rename the 4th sheet from "OrderClientRef" to "OrderReference"
$wsRef = $wb.Worksheets.Item("OrderClientRef")
$wsRef.Name = "OrderReference"

# Add new shared string / header "TotalExpense" into C1
$wsRef.Range("C1").Value = "TotalExpense"

# Fill in sequence values A2:A5 = 1,2,3,4
$wsRef.Range("A2").Value = 1
$wsRef.Range("A3").Value = 2
$wsRef.Range("A4").Value = 3
$wsRef.Range("A5").Value = 4

# Update selection on this sheet
$wsRef.Range("F15").Select()

# Make this sheet the active/selected tab
$wsRef.Activate()

$wb.Save()
